# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 14 of the data table
# (pushing the existing rows 14-43 down to rows 15-44), and the new row
# is populated with the latest week's Ciboulette price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14; everything below (rows 14-43) shifts down
# to rows 15-44, carrying its formatting (including the date number
# format on column D) along with it.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new week's record.
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C14").Value = "Ñuble"
$ws.Range("D14").Value = 45251
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = 100112039
$ws.Range("G14").Value = "Ciboulette"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("N14").Value = "`$/docena de atados"
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 833
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = "Hortaliza"
